# temps_labos.xlsx — "deploy labo 1,2,3 on the shinyapp.io of uqar"
#
# Log two new work-log entries on 2024-08-26 (serial 45530):
#   row 27: 5h  "Révision labo 3 et ajustement document labo 4"
#   row 28: 3h  "Rencontre Joël et deploiement shiny"
#
# Row 28 is the new most-recent entry, so it keeps the distinct "latest
# row" date style that used to sit on rows 25-26; those fall back to the
# ordinary date style shared by the rest of column A, as does the other
# newly typed date in row 27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting two rows below the last existing row (instead of just typing
# into previously-empty cells) makes Excel copy row 26's formatting
# (including its distinctive date style) down into the new rows 27-28 -
# that's what lets row 28 end up the "latest" highlighted row below.
$ws.Rows("27:28").Insert()

$ws.Range("A27").Value = 45530
$ws.Range("B27").Value = 5
$ws.Range("C27").Value = "Révision labo 3 et ajustement document labo 4"

$ws.Range("A28").Value = 45530
$ws.Range("B28").Value = 3
$ws.Range("C28").Value = "Rencontre Joël et deploiement shiny"

# Rows 25, 26 and the new row 27 are no longer the latest entry, so they
# go back to the plain date format used everywhere else in column A.
$ws.Range("A25").NumberFormat = "yyyy/mm/dd"
$ws.Range("A26").NumberFormat = "yyyy/mm/dd"
$ws.Range("A27").NumberFormat = "yyyy/mm/dd"

# Inserting rows widened the SUM range B2:B200 that already safely covers
# rows 27-28; put it back to the original formula text (same result).
$ws.Range("D2").Formula = "=SUM(B2:B200)"

# Leave the cursor where the author finished typing.
[void]$ws.Range("C28").Select()
